# Add a new row of data (row 3) with two attorney/email pairs,
# the emails being mailto hyperlinks (matching Excel's auto-hyperlink + Hyperlink style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "Juan Perez"
$ws.Range("F3").Value = "juan.perez@gmail.com"
$null = $ws.Hyperlinks.Add($ws.Range("F3"), "mailto:juan.perez@gmail.com", "", "", "juan.perez@gmail.com")

$ws.Range("J3").Value = "Jose Garcia"
$ws.Range("K3").Value = "jose.garcia@gmail.com"
$null = $ws.Hyperlinks.Add($ws.Range("K3"), "mailto:jose.garcia@gmail.com", "", "", "jose.garcia@gmail.com")

$null = $ws.Range("A3").Select()
